$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Write column-by-column so new shared strings are interned in the same
# order as the target workbook (A23, A24, B23, B24, C23, C24, D23, D24).
$ws.Range("A23").Value = "PUBLONS023"
$ws.Range("A24").Value = "PUBLONS024"

$ws.Range("B23").Value = "OPQA-5891"
$ws.Range("B24").Value = "OPQA-5892"

$ws.Range("C23").Value = "Verify Error message when user  exist in platform in Blocked state and trying to login."
$ws.Range("C24").Value = "Verify Error message when user  exist in platform in Evicted state and trying to login."

$ws.Range("D23").Value = "Y"
$ws.Range("D24").Value = "Y"

# Update the selection to cover D22:D24 with the active cell at D22
$ws.Range("D22:D24").Select()
